# "Generate Report for Handback" - update localization-status report rows
# for the file 59e2ffa5-d6fd-4842-a146-a4f49be793e5.md now that the
# handback has completed (no longer just "Ready for handoff").

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = "Handed back: in sync with en-US"
$ovw.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("K3").Value = "2016-09-05 16:57:32"
$zh.Range("P3").Value = ""

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("K3").Value = "2016-09-05 16:57:39"
$de.Range("P3").Value = ""
